# Updates cryptos list figures (price + 1h volume change) to match the latest
# scrape, and re-orders the final three coin rows (ordi / FraxShare / EnergySwap
# -> FraxShare / ordi / MultiversX).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "43.979.28"
$ws.Range("E2").Value = "  -0.19%  "
# Row 3
$ws.Range("D3").Value = "2.234.57"
$ws.Range("E3").Value = "  -0.57%  "
# Row 4
$ws.Range("E4").Value = "  +0.23%  "
# Row 5
$ws.Range("D5").Formula = "'306.06"
$ws.Range("E5").Value = "  -3.85%  "
# Row 6
$ws.Range("D6").Formula = "'94.46"
$ws.Range("E6").Value = "  -6.34%  "
# Row 7
$ws.Range("D7").Formula = "'0.570"
$ws.Range("E7").Value = "  -0.61%  "
# Row 8
$ws.Range("E8").Value = "  +0.23%  "
# Row 9
$ws.Range("D9").Formula = "'0.522"
$ws.Range("E9").Value = "  -3.95%  "
# Row 10
$ws.Range("D10").Formula = "'34.71"
$ws.Range("E10").Value = "  -5.61%  "
# Row 11
$ws.Range("D11").Formula = "'0.0808"
$ws.Range("E11").Value = "  -2.22%  "
# Row 12
$ws.Range("D12").Formula = "'7.19"
$ws.Range("E12").Value = "  -4.15%  "
# Row 13
$ws.Range("E13").Value = "  -0.79%  "
# Row 14
$ws.Range("D14").Value = "2.576.84"
$ws.Range("E14").Value = "  -0.58%  "
# Row 15
$ws.Range("D15").Value = "2.235.86"
$ws.Range("E15").Value = "  -2.23%  "
# Row 16
$ws.Range("D16").Formula = "'0.821"
$ws.Range("E16").Value = "  -3.20%  "
# Row 17
$ws.Range("D17").Formula = "'13.53"
$ws.Range("E17").Value = "  -4.65%  "
# Row 18
$ws.Range("D18").Value = "43.866.42"
$ws.Range("E18").Value = "  -0.22%  "
# Row 19
$ws.Range("D19").Value = "0.0₃0961"
$ws.Range("E19").Value = "  -1.32%  "
# Row 20
$ws.Range("D20").Formula = "'11.99"
$ws.Range("E20").Value = "  -11.05%  "
# Row 21
$ws.Range("D21").Formula = "'6.26"
$ws.Range("E21").Value = "  -2.89%  "
# Row 22
$ws.Range("D22").Formula = "'65.17"
$ws.Range("E22").Value = "  -0.30%  "
# Row 23
$ws.Range("D23").Formula = "'236.49"
$ws.Range("E23").Value = "  +0.78%  "
# Row 24
$ws.Range("E24").Value = "  -5.34%  "
# Row 25
$ws.Range("E25").Value = "  -4.73%  "
# Row 26
$ws.Range("E26").Value = "  +0.23%  "
# Row 27
$ws.Range("D27").Formula = "'9.87"
$ws.Range("E27").Value = "  -5.88%  "
# Row 28
$ws.Range("D28").Formula = "'37.70"
$ws.Range("E28").Value = "  -0.97%  "
# Row 29
$ws.Range("E29").Value = "  -2.05%  "
# Row 30
$ws.Range("D30").Formula = "'5.97"
$ws.Range("E30").Value = "  -1.54%  "
# Row 31
$ws.Range("D31").Formula = "'19.90"
$ws.Range("E31").Value = "  -0.90%  "
# Row 32
$ws.Range("D32").Formula = "'152.43"
$ws.Range("E32").Value = "  -3.81%  "
# Row 33
$ws.Range("D33").Formula = "'0.0800"
$ws.Range("E33").Value = "  -5.44%  "
# Row 34
$ws.Range("D34").Formula = "'3.24"
$ws.Range("E34").Value = "  +0.78%  "
# Row 35
$ws.Range("D35").Formula = "'2.59"
$ws.Range("E35").Value = "  -3.26%  "
# Row 36
$ws.Range("D36").Formula = "'0.109"
$ws.Range("E36").Value = "  -3.68%  "
# Row 37
$ws.Range("E37").Value = "  +0.59%  "
# Row 38
$ws.Range("D38").Formula = "'1.77"
$ws.Range("E38").Value = "  -8.98%  "
# Row 39
$ws.Range("D39").Formula = "'14.94"
$ws.Range("E39").Value = "  -6.89%  "
# Row 40
$ws.Range("D40").Formula = "'3.81"
$ws.Range("E40").Value = "  -7.87%  "
# Row 41
$ws.Range("E41").Value = "  -9.07%  "
# Row 42
$ws.Range("D42").Formula = "'0.0298"
$ws.Range("E42").Value = "  -4.90%  "
# Row 43
$ws.Range("E43").Value = "  +0.33%  "
# Row 44
$ws.Range("D44").Value = "1.720.69"
$ws.Range("E44").Value = "  -2.18%  "
# Row 45
$ws.Range("D45").Formula = "'84.58"
$ws.Range("E45").Value = "  +4.19%  "
# Row 46
$ws.Range("D46").Formula = "'0.187"
$ws.Range("E46").Value = "  -4.40%  "
# Row 47
$ws.Range("D47").Formula = "'99.61"
$ws.Range("E47").Value = "  -3.19%  "
# Row 48
$ws.Range("D48").Formula = "'4.93"
$ws.Range("E48").Value = "  -4.30%  "
# Row 49
$ws.Range("B49").Value = "FraxShare"
$ws.Range("C49").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D49").Formula = "'8.03"
$ws.Range("E49").Value = "  -3.19%  "
# Row 50
$ws.Range("B50").Value = "ordi"
$ws.Range("C50").Value = "https://coinranking.com/coin/j7-7vPrOi+ordi-ordi"
$ws.Range("D50").Formula = "'68.44"
$ws.Range("E50").Value = "  -8.16%  "
# Row 51
$ws.Range("B51").Value = "MultiversX"
$ws.Range("C51").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D51").Formula = "'54.01"
$ws.Range("E51").Value = "  -5.81%  "
